# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Every cell in C2:C436 held the serial date value 46060 and is bumped
# to 46061 (i.e. the "last updated" timestamp advanced by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 436 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
